# "Generate Report for Archive" — refresh the localization-status report:
# flip any "Ready for handoff" status cell to "In Translation" across every
# sheet (Overview + each locale tab), then re-fit the Status columns that
# held the old, wider text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ("Ready for handoff" -eq $cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# The status text got shorter, so re-fit the columns that used to size
# themselves around "Ready for handoff".
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns("C:C").ColumnWidth = 12.5
